# update scripts wuth new tpm
# Refresh the NATMI ligand/receptor metrics on the sole data sheet with
# values recomputed from the new TPM expression matrix.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs -> ECs)
$ws.Range("G2").Value = 0.158418
$ws.Range("H2").Value = 0.475254
$ws.Range("I2").Value = 0.05652797120826585
$ws.Range("J2").Value = 0.05652797120826585
$ws.Range("M2").Value = 0.4890553333333333
$ws.Range("N2").Value = 1.467166
$ws.Range("O2").Value = 0.9644476581758422
$ws.Range("P2").Value = 0.9644476581758422
$ws.Range("Q2").Value = 0.07747516779599999
$ws.Range("R2").Value = 0.6972765101639999
$ws.Range("S2").Value = 0.05451826945324343
$ws.Range("T2").Value = 0.05451826945324343

# Row 3 (ECs -> MuSCs)
$ws.Range("G3").Value = 0.158418
$ws.Range("H3").Value = 0.475254
$ws.Range("I3").Value = 0.05652797120826585
$ws.Range("J3").Value = 0.05652797120826585
$ws.Range("O3").Value = 0.03555234182415776
$ws.Range("P3").Value = 0.03555234182415776
$ws.Range("Q3").Value = 0.002855959703999999
$ws.Range("R3").Value = 0.025703637336
$ws.Range("S3").Value = 0.002009701755022416
$ws.Range("T3").Value = 0.002009701755022416

# Row 4 (FAPs -> ECs)
$ws.Range("I4").Value = 0.7628354881578912
$ws.Range("J4").Value = 0.7628354881578912
$ws.Range("M4").Value = 0.4890553333333333
$ws.Range("N4").Value = 1.467166
$ws.Range("O4").Value = 0.9644476581758422
$ws.Range("P4").Value = 0.9644476581758422
$ws.Range("Q4").Value = 1.045514391946444
$ws.Range("R4").Value = 9.409629527518
$ws.Range("S4").Value = 0.7357149001273036
$ws.Range("T4").Value = 0.7357149001273036

# Row 5 (FAPs -> MuSCs)
$ws.Range("I5").Value = 0.7628354881578912
$ws.Range("J5").Value = 0.7628354881578912
$ws.Range("O5").Value = 0.03555234182415776
$ws.Range("P5").Value = 0.03555234182415776
$ws.Range("S5").Value = 0.0271205880305876
$ws.Range("T5").Value = 0.0271205880305876

# Row 6 (MuSCs -> ECs)
$ws.Range("G6").Value = 0.5062286666666667
$ws.Range("I6").Value = 0.180636540633843
$ws.Range("J6").Value = 0.180636540633843
$ws.Range("M6").Value = 0.4890553333333333
$ws.Range("N6").Value = 1.467166
$ws.Range("O6").Value = 0.9644476581758422
$ws.Range("P6").Value = 0.9644476581758422
$ws.Range("Q6").Value = 0.2475738293195555
$ws.Range("R6").Value = 2.228164463876
$ws.Range("S6").Value = 0.1742144885952953
$ws.Range("T6").Value = 0.1742144885952953

# Row 7 (MuSCs -> MuSCs)
$ws.Range("G7").Value = 0.5062286666666667
$ws.Range("I7").Value = 0.180636540633843
$ws.Range("J7").Value = 0.180636540633843
$ws.Range("O7").Value = 0.03555234182415776
$ws.Range("P7").Value = 0.03555234182415776
$ws.Range("Q7").Value = 0.009126290402666666
$ws.Range("R7").Value = 0.082136613624
$ws.Range("S7").Value = 0.00642205203854775
$ws.Range("T7").Value = 0.00642205203854775
